{"js": "// Update the cover letter for a new submission:\n// 1. Change the target journal name from \"Series A: Statistics in Society\"\n//    to \"Series C: Applied Statistics\" (keeps the existing italic formatting).\n// 2. Drop the trailing \", published in your journal\" clause about the 2018\n//    paper reference, since the letter is now addressed to a different journal.\n\nconst body = context.document.body;\n\n// 1) Retarget the journal name (italic run).\nconst journalHits = body.search(\"Series A: Statistics in Society\", { matchCase: true });\njournalHits.load(\"items\");\nawait context.sync();\n\nfor (const hit of journalHits.items) {\n  hit.insertText(\"Series C: Applied Statistics\", \"Replace\");\n}\nawait context.sync();\n\n// 2) Remove \", published in your journal\" after the 2018 paper reference.\nconst publishedHits = body.search(\", published in your journal\", { matchCase: true });\npublishedHits.load(\"items\");\nawait context.sync();\n\nfor (const hit of publishedHits.items) {\n  hit.insertText(\"\", \"Replace\");\n}\nawait context.sync();\n", "ps1": "# Update the cover letter for a new submission:\n# 1. Change the target journal name from \"Series A: Statistics in Society\"\n#    to \"Series C: Applied Statistics\" (keeps the existing italic formatting).\n# 2. Drop the trailing \", published in your journal\" clause about the 2018\n#    paper reference, since the letter is now addressed to a different journal.\n\n$d = $word.ActiveDocument\n\n# 1) Retarget the journal name (italic run).\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Execute(\"Series A: Statistics in Society\", $false, $false, $false, $false, $false, $true, 1, $false, \"Series C: Applied Statistics\", 2)\n\n# 2) Remove \", published in your journal\" after the 2018 paper reference.\n$find2 = $d.Content.Find\n$find2.ClearFormatting()\n$find2.Replacement.ClearFormatting()\n$find2.Execute(\", published in your journal\", $false, $false, $false, $false, $false, $true, 1, $false, \"\", 2)\n"}
